# Weekly update: a new "Ajo" (Chino/Primera) price record for
# Terminal La Palmera de La Serena is inserted as row 155, pushing the
# existing records (old rows 155-191) down by one row (new rows 156-192).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 155; this shifts rows
# 155..191 down to 156..192 and keeps all of their values/styles intact.
$ws.Rows.Item(155).Insert()

# Populate the new row 155 with the new weekly record.
$ws.Range("A155").Value = 8
$ws.Range("B155").Value = "Terminal La Palmera de La Serena"
$ws.Range("C155").Value = "Coquimbo"
$ws.Range("D155").Value = 44551
$ws.Range("E155").Value = 4
$ws.Range("F155").Value = 100112003
$ws.Range("G155").Value = "Ajo"
$ws.Range("H155").Value = "Chino"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 560
$ws.Range("K155").Value = 19000
$ws.Range("L155").Value = 19500
$ws.Range("M155").Value = 19250
$ws.Range("N155").Value = "`$/caja 10 kilos"
$ws.Range("O155").Value = "China"
$ws.Range("P155").Value = 1925
$ws.Range("Q155").Value = 10
$ws.Range("R155").Value = "Hortaliza"
